$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lom3229 = "LOM3229 -  Métodos Experimentais da Física II  (Indicação de Conjunto)`n"
$lob1021 = "LOB1021 -  Física IV  (Requisito)`n"
$lom3016 = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"

$ws.Range("B23").Value = $lom3229
$ws.Range("C23").Value = $lom3229

$ws.Range("B24").Value = $lob1021
$ws.Range("C24").Value = $lob1021

$ws.Range("B25").Value = $lom3016
$ws.Range("C25").Value = $lom3016
